# Update computed result values in the "pl_mw" sheet (res_line) for the
# "case with 380 kV" scenario. Only the numeric result cells (columns
# B,C,D,E,G,I,K,N,O for rows 2-25) change; other cells are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3728393510611738
$ws.Range("C2").Value = 0.07131736678489631
$ws.Range("D2").Value = 0.07872640225792793
$ws.Range("E2").Value = 0.4147455181249313
$ws.Range("G2").Value = 0.002430260566710424
$ws.Range("I2").Value = 0.5652666147676086
$ws.Range("K2").Value = 0.4167048701545468
$ws.Range("N2").Value = 1.354487868334006
$ws.Range("O2").Value = 2.722997057814212
$ws.Range("B3").Value = 0.3317514051404373
$ws.Range("C3").Value = 0.06216950224509787
$ws.Range("D3").Value = 0.07138503577733957
$ws.Range("E3").Value = 0.3618866571516293
$ws.Range("G3").Value = 0.002433074920849654
$ws.Range("I3").Value = 0.5680074896709328
$ws.Range("K3").Value = 0.368047700144956
$ws.Range("N3").Value = 1.37095577368613
$ws.Range("O3").Value = 2.719235485481647
$ws.Range("B4").Value = 0.3065718528187631
$ws.Range("C4").Value = 0.0565438193407175
$ws.Range("D4").Value = 0.06691191868283397
$ws.Range("E4").Value = 0.3295189922669834
$ws.Range("G4").Value = 0.002434893890683622
$ws.Range("I4").Value = 0.5700182368928424
$ws.Range("K4").Value = 0.3382027217104167
$ws.Range("N4").Value = 1.381572108247143
$ws.Range("O4").Value = 2.718637902293807
$ws.Range("B5").Value = 0.2963235066082177
$ws.Range("C5").Value = 0.05424905112680278
$ws.Range("D5").Value = 0.06509775574764376
$ws.Range("E5").Value = 0.3163493327397617
$ws.Range("G5").Value = 0.002435658077610145
$ws.Range("I5").Value = 0.5709199546523713
$ws.Range("K5").Value = 0.3260486082881471
$ws.Range("N5").Value = 1.386025226058939
$ws.Range("O5").Value = 2.718824129305858
$ws.Range("B6").Value = 0.2946225420151052
$ws.Range("C6").Value = 0.0538678694727821
$ws.Range("D6").Value = 0.06479703866743591
$ws.Range("E6").Value = 0.3141637057907332
$ws.Range("G6").Value = 0.002435786357834863
$ws.Range("I6").Value = 0.5710746528728521
$ws.Range("K6").Value = 0.3240309145883487
$ws.Range("N6").Value = 1.386772322823211
$ws.Range("O6").Value = 2.718880985766219
$ws.Range("B7").Value = 0.3064335889200152
$ws.Range("C7").Value = 0.05651288046932734
$ws.Range("D7").Value = 0.06688741716055802
$ws.Range("E7").Value = 0.3293413011213318
$ws.Range("G7").Value = 0.002434904103945769
$ws.Range("I7").Value = 0.5700300646114087
$ws.Range("K7").Value = 0.338038774555514
$ws.Range("N7").Value = 1.381631650981044
$ws.Range("O7").Value = 2.718638674804794
$ws.Range("B8").Value = 0.3586623054264066
$ws.Range("C8").Value = 0.06816500071714415
$ws.Range("D8").Value = 0.07618792805907049
$ws.Range("E8").Value = 0.3965005171208702
$ws.Range("G8").Value = 0.00243121212426556
$ws.Range("I8").Value = 0.5661435547562412
$ws.Range("K8").Value = 0.3999216257291494
$ws.Range("N8").Value = 1.360061133238678
$ws.Range("O8").Value = 2.72134426701399
$ws.Range("B9").Value = 0.461460398411873
$ws.Range("C9").Value = 0.09094696628088172
$ws.Range("D9").Value = 0.09470146326171403
$ws.Range("E9").Value = 0.5289785444732757
$ws.Range("G9").Value = 0.00242469044018236
$ws.Range("I9").Value = 0.561129033643688
$ws.Range("K9").Value = 0.5215134720082517
$ws.Range("N9").Value = 1.321772167710814
$ws.Range("O9").Value = 2.740272853758114
$ws.Range("B10").Value = 0.5372135153181716
$ws.Range("C10").Value = 0.1076485985971374
$ws.Range("D10").Value = 0.108474572119718
$ws.Range("E10").Value = 0.6269069576437971
$ws.Range("G10").Value = 0.002420332162151052
$ws.Range("I10").Value = 0.5590425688746734
$ws.Range("K10").Value = 0.6109970235971502
$ws.Range("N10").Value = 1.296090868192101
$ws.Range("O10").Value = 2.762544573741337
$ws.Range("B11").Value = 0.5717248311590311
$ws.Range("C11").Value = 0.1152399216067295
$ws.Range("D11").Value = 0.1147783092507524
$ws.Range("E11").Value = 0.671613665167655
$ws.Range("G11").Value = 0.002418442537731786
$ws.Range("I11").Value = 0.5584421937062629
$ws.Range("K11").Value = 0.6517394140654176
$ws.Range("N11").Value = 1.284940684881857
$ws.Range("O11").Value = 2.774506357543402
$ws.Range("B12").Value = 0.5848004943735248
$ws.Range("C12").Value = 0.1181137083724195
$ws.Range("D12").Value = 0.1171709084847237
$ws.Range("E12").Value = 0.688567746433705
$ws.Range("G12").Value = 0.002417740280490739
$ws.Range("I12").Value = 0.558265146551733
$ws.Range("K12").Value = 0.6671725667901853
$ws.Range("N12").Value = 1.280795078677814
$ws.Range("O12").Value = 2.779300151853676
$ws.Range("B13").Value = 0.5819841113455766
$ws.Range("C13").Value = 0.1174948257667552
$ws.Range("D13").Value = 0.1166553742415886
$ws.Range("E13").Value = 0.6849152575246933
$ws.Range("G13").Value = 0.00241789093358214
$ws.Range("I13").Value = 0.5583010372403692
$ws.Range("K13").Value = 0.6638485477581071
$ws.Range("N13").Value = 1.281684493944097
$ws.Range("O13").Value = 2.778255960390368
$ws.Range("B14").Value = 0.5728004359205556
$ws.Range("C14").Value = 0.1154763672975889
$ws.Range("D14").Value = 0.1149750393509379
$ws.Range("E14").Value = 0.6730079837155643
$ws.Range("G14").Value = 0.002418384496406541
$ws.Range("I14").Value = 0.5584266189452052
$ws.Range("K14").Value = 0.6530090122597301
$ws.Range("N14").Value = 1.284598084656619
$ws.Range("O14").Value = 2.774895447011261
$ws.Range("B15").Value = 0.5671760682638194
$ws.Range("C15").Value = 0.1142398891256278
$ws.Range("D15").Value = 0.1139465035557805
$ws.Range("E15").Value = 0.6657176953574435
$ws.Range("G15").Value = 0.002418688548277278
$ws.Range("I15").Value = 0.5585100966127854
$ws.Range("K15").Value = 0.6463701132663289
$ws.Range("N15").Value = 1.28639274029176
$ws.Range("O15").Value = 2.772871460501051
$ws.Range("B16").Value = 0.5349591141335281
$ws.Range("C16").Value = 0.1071523627361728
$ws.Range("D16").Value = 0.1080633789026422
$ws.Range("E16").Value = 0.623988615222089
$ws.Range("G16").Value = 0.002420457518663772
$ws.Range("I16").Value = 0.5590888356586916
$ws.Range("K16").Value = 0.6083351117285645
$ws.Range("N16").Value = 1.296830286151359
$ws.Range("O16").Value = 2.761799752544533
$ws.Range("B17").Value = 0.515207839883459
$ws.Range("C17").Value = 0.1028027925753463
$ws.Range("D17").Value = 0.1044640814753706
$ws.Range("E17").Value = 0.5984310806177433
$ws.Range("G17").Value = 0.002421566490594598
$ws.Range("I17").Value = 0.559533307168735
$ws.Range("K17").Value = 0.5850108690580385
$ws.Range("N17").Value = 1.303369856229684
$ws.Range("O17").Value = 2.755477057341864
$ws.Range("B18").Value = 0.5038522347686296
$ws.Range("C18").Value = 0.1003004370052452
$ws.Range("D18").Value = 0.1023974636158442
$ws.Range("E18").Value = 0.5837458794870827
$ws.Range("G18").Value = 0.00242221309749114
$ws.Range("I18").Value = 0.5598217741106808
$ws.Range("K18").Value = 0.571598773884233
$ws.Range("N18").Value = 1.307181342470526
$ws.Range("O18").Value = 2.752012627155779
$ws.Range("B19").Value = 0.5000082575502915
$ws.Range("C19").Value = 0.09945307887755916
$ws.Range("D19").Value = 0.1016983609186894
$ws.Range("E19").Value = 0.5787762123800206
$ws.Range("G19").Value = 0.002422433533290901
$ws.Range("I19").Value = 0.5599250761711119
$ws.Range("K19").Value = 0.5670582596586371
$ws.Range("N19").Value = 1.308480445507421
$ws.Range("O19").Value = 2.750869180650454
$ws.Range("B20").Value = 0.5173099018293215
$ws.Range("C20").Value = 0.1032658734176266
$ws.Range("D20").Value = 0.1048468599084202
$ws.Range("E20").Value = 0.6011501772280639
$ws.Range("G20").Value = 0.002421447532940808
$ws.Range("I20").Value = 0.5594825947249618
$ws.Range("K20").Value = 0.5874934253464232
$ws.Range("N20").Value = 1.302668521581777
$ws.Range("O20").Value = 2.756132288466517
$ws.Range("B21").Value = 0.575497717310725
$ws.Range("C21").Value = 0.1160692611994705
$ws.Range("D21").Value = 0.1154684448308956
$ws.Range("E21").Value = 0.6765047558845225
$ws.Range("G21").Value = 0.002418239164507572
$ws.Range("I21").Value = 0.5583883661781357
$ws.Range("K21").Value = 0.6561927158747949
$ws.Range("N21").Value = 1.283740208718374
$ws.Range("O21").Value = 2.775875335435956
$ws.Range("B22").Value = 0.6135673315839085
$ws.Range("C22").Value = 0.1244318850102957
$ws.Range("D22").Value = 0.1224423884165162
$ws.Range("E22").Value = 0.7258976970347391
$ws.Range("G22").Value = 0.00241621981751683
$ws.Range("I22").Value = 0.557966492719828
$ws.Range("K22").Value = 0.7011200767151706
$ws.Range("N22").Value = 1.271816763821906
$ws.Range("O22").Value = 2.790318528637158
$ws.Range("B23").Value = 0.5932452649031177
$ws.Range("C23").Value = 0.1199690551808033
$ws.Range("D23").Value = 0.1187173221769058
$ws.Range("E23").Value = 0.6995219401091077
$ws.Range("G23").Value = 0.002417290511325043
$ws.Range("I23").Value = 0.5581647707735158
$ws.Range("K23").Value = 0.6771389866478614
$ws.Range("N23").Value = 1.278139544738686
$ws.Range("O23").Value = 2.782468705863863
$ws.Range("B24").Value = 0.5163595606414333
$ws.Range("C24").Value = 0.1030565199603757
$ws.Range("D24").Value = 0.1046737974763232
$ws.Range("E24").Value = 0.5999208483873986
$ws.Range("G24").Value = 0.002421501285501004
$ws.Range("I24").Value = 0.5595054192383415
$ws.Range("K24").Value = 0.5863710699707383
$ws.Range("N24").Value = 1.30298543362264
$ws.Range("O24").Value = 2.755835527199451
$ws.Range("B25").Value = 0.4336106611200705
$ws.Range("C25").Value = 0.08479057997772088
$ws.Range("D25").Value = 0.08966324012212112
$ws.Range("E25").Value = 0.4930437835459998
$ws.Range("G25").Value = 0.002426378317122093
$ws.Range("I25").Value = 0.5622056366039274
$ws.Range("K25").Value = 0.4885935989583459
$ws.Range("N25").Value = 1.331700513827861
$ws.Range("O25").Value = 2.733687568917674
